# Insert a new weekly record at row 162 of Sheet1, pushing existing rows
# 162-195 down to 163-196 (all their data/formatting shifts automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 162; Excel shifts rows 162:195 down to 163:196.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new record's values.
$ws.Cells.Item(162, 1).Value = 7
$ws.Cells.Item(162, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(162, 3).Value = 'Ñuble'
$ws.Cells.Item(162, 4).Value = 45204
$ws.Cells.Item(162, 5).Value = 16
$ws.Cells.Item(162, 6).Value = 'Fruta'
$ws.Cells.Item(162, 7).Value = 100108
$ws.Cells.Item(162, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(162, 9).Value = 100108002
$ws.Cells.Item(162, 10).Value = 'Mango'
$ws.Cells.Item(162, 11).Value = 'Sin especificar'
$ws.Cells.Item(162, 12).Value = 'Primera'
$ws.Cells.Item(162, 13).Value = 30
$ws.Cells.Item(162, 14).Value = 10000
$ws.Cells.Item(162, 15).Value = 10000
$ws.Cells.Item(162, 16).Value = 10000
$ws.Cells.Item(162, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(162, 18).Value = 'Brasil'
$ws.Cells.Item(162, 19).Value = 2500
$ws.Cells.Item(162, 20).Value = 4
